# Applies the SSS artefact edits described in the commit:
# "Modificando artefato 22, corrigindo as regras de acordo com 2 cenarios"
#
#   SSS-00001: troca a regra de "permitir compra" para "descartar
#              solicitacao de orcamento por inatividade".
#   SSS-00002: troca a regra de "consultar produtos disponiveis" para
#              "loja gerar orcamento para consulta previa".
#   SSS-00003: troca "o cliente" por "a loja" como quem cancela o pedido.

$d = $word.ActiveDocument

$wdFindContinue   = 1
$wdReplaceAll     = 2

# --- SSS-00001 --------------------------------------------------------
$d.Content.Find.Execute(
    "permitir que o cliente gere um orçamento antes de definitivamente efetuar a compra.",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "descartar a solicitação de orçamento dentro de 1h de inatividade do sistema.",
    $wdReplaceAll)

# --- SSS-00002 --------------------------------------------------------
$d.Content.Find.Execute(
    "consultar os produtos que estão disponíveis.",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "permitir que a loja gere um orçamento para consulta previa, antes de ser aprovado.",
    $wdReplaceAll)

# --- SSS-00003 --------------------------------------------------------
$d.Content.Find.Execute(
    "que o cliente cancele o seu pedido.",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "que a loja cancele o seu pedido.",
    $wdReplaceAll)
